$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force the cell to store a literal text value even when the text looks
    # like a number (e.g. "1.1", "-3.1"), mirroring how these generator
    # workbooks keep numeric-looking expressions as plain strings. We flip
    # the cell to Text format just long enough for the assignment to stick,
    # then clear the formatting back to the workbook default so no visible
    # style change remains on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# --- Sheet: Restricciones_del_lider ---
$ws2 = $wb.Worksheets.Item("Restricciones_del_lider")
$ws2.Range("A2").Value = "-2.1 + x"
Set-TextValue $ws2.Range("B2") "1.1"
Set-TextValue $ws2.Range("D2") "0.74"
$ws2.Range("A3").Value = "2.1 - x"
Set-TextValue $ws2.Range("B3") "-3.1"
Set-TextValue $ws2.Range("D3") "0.27"

# --- Sheet: Restricciones_del_follower ---
$ws3 = $wb.Worksheets.Item("Restricciones_del_follower")
$ws3.Range("A2").Value = "-10.025316455696203 + 3.0379746835443036y"
Set-TextValue $ws3.Range("B2") "9.025316455696203"
Set-TextValue $ws3.Range("D2") "0.22"
Set-TextValue $ws3.Range("E2") "3.8"
Set-TextValue $ws3.Range("F2") "7.199999999999999"
$ws3.Range("A3").Value = "-4.521000000000001 + 1.37y"
Set-TextValue $ws3.Range("B3") "3.521000000000001"
Set-TextValue $ws3.Range("D3") "0.66"
Set-TextValue $ws3.Range("E3") "6.0"
Set-TextValue $ws3.Range("F3") "4.4"

# --- Sheet: Punto_modificado ---
$ws4 = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws4.Range("A2") "2.1"
Set-TextValue $ws4.Range("B2") "3.3000000000000003"

# --- Sheet: Vector_bf ---
# NOTE: "Vector_bf" and "Vector_BF" only differ by case, and
# Worksheets.Item(<name>) resolves case-insensitively, so these two sheets
# are addressed by their 1-based tab position instead of by name.
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "-10.362554430379749"

# --- Sheet: Vector_BF ---
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "-1.7999999999999998"
Set-TextValue $ws6.Range("A3") "-42.006303797468355"

# --- Sheet: Vector_Alpha ---
$ws7 = $wb.Worksheets.Item("Vector_Alpha")
$ws7.Range("A2").Value = 2.37
